# =====================================================================
# Edit script: applies the 'research3 exercise diet' Jan/Feb 2021 commit
# =====================================================================
$wb = $excel.ActiveWorkbook

$wsNutrition = $wb.Worksheets.Item("NutritionalData")
$wsMeasures  = $wb.Worksheets.Item("researchMeasures")

# --- researchMeasures: update existing row 71 text -----------------
$zTextRow71 = @'
Woke up 610 am got out of bed to feed the babies and put Growly's heart meds in his food at the same time everyday. Alarm went off at 6:30 am, went back to bed by about 645 am and got out of bed at 830 am. I went to bed at 330 am after spending 2 hours designing some window decals to promote my side biz and also try out new marketing ideas to grab attention. So I probably got about 5 hours but more like 4 1/2 hours of sleep. Got up and made coffee and folded laundry from yesterday's linens. Fed the cat outside because I saw her outside. She has a limp, poor baby, but has had it for a few days now. Noticed it Tuesday. Still a little spotty, not light but spotty. Didn't have a lg BM until 11 am approx. while drinking 3rd cup of coffee and after having a large bowl of the last of the pasta made a few days ago on 3-23-21. Started first half of week 6 ppt in genetics Bio-18. Tired now by 1130 am, and have all day of class starting at 3 pm. Might take a nap later. Measurements around 1 pm. Ate 2 quesadillas with cinnamon and paprika at around 1215 pm and shared last half of choco banana cake with babies. Finished the 2nd part of genetics recorded lecture. Made another pot of pasta with gluten free spaghetti, prego 3 cheese, 2 bell peppers orange and red, and 1 pkg beyond meat at break time between lecture and lab. The plumber David came by earlier before 1 pm to install a water heater, but Dave didn't tell us anything about it. He left the old one outside or the same one. But I heard him drilling and making banging noise under the house. I had made some cold pressed coffee earlier in the morning before starting the lectures in genetics that I drank a cup after taking a 10 minute nap around 1 pm. I put some organic sugar I had about a year now in the cabinet and some sour cream in it. I know, but no dairy. It didn't dissolve, so not an aqueous solution. But tasted good. I finally took out the kitchen stuff and mat from the IKEA trip last week early in the day into use. Only washed the coffee mug. I don't trust putting hot water into it, because it could be pyrite and not sure, but could break or crack. Like a pyrite bowl I put in the microwave last time I bought kitchen stuff there. I also took out all the note cards, pencils, pens, erasers, sketch pad, etc from Staples that was still in the bag a week, and the two decorative boxes from Michael's craft store bought at same time.  I was able to copy my notes on the videos of chemical reactions from my notebook into my personal lab manual before lecture and put together the facial machine, but screw up the wheel. Put it in my car and fixed the wheel with my personal tools and secured it to the inside. Its a light weight version of the pro quality facial steamers at work. I think I got my hydrocollator that I spent $200 on and it looks really tiny, smaller than a shoe box. This should be interesting. But it got great reviews, so could work fine for single or couples to plug it in and use the heating pads on their backs. I will have to upload a photo. I just opened it and forgot that I was also expecting a wifi booster for the roommate's room since he bitches and complains and gets on my case about getting it fixed. They sent me the wifi. Thats great, because that would have been a super tiny hydrocollator if it was. I have ordered two other items like my kabrow balm for eyebrows from benefit that was tinier and the aquagel for the RF machine that I thought would be the size of a regular 16 oz bottle of rubbing alcohol, but was more like a hand sanitizer portable size. It was misleading because it said 64 g, and that was confused with 64 oz. Not the same thing. I only saved $3-4 for ordering the tiny versions. Those items came in about 2 weeks ago. Bed time around 930 pm shortly after class ended. Mom called 10 minutes before class got out. 
'@
$wsMeasures.Range("Z71").Value2 = $zTextRow71

$aaTextRow71 = @'
1 1/2 bowls pasta
(1152.00	39.28	10.03	73.31	137.44	25.31	1426.88)
4 tbs sourcream
(120	10	7	2	4	0	30)
4 corn tortillas Guerrero
(200	2	0	4	42	4	40)
1/3 cup mozzarella cheese
(106.67	6.67	4.67	8.00	1.33	0.00	253.33)
1/2 chocobanana cake from 85degrees celsius
(330	14	11	5	44	0	60)

1 cup sirens starbucks medium roast French ground coffee cold pressed 5 hours in fridge
with 1 tbs pure can sugar full circle brand
pure can sugar
(45	0	0	0	12	0	0)
1 tbs sour cream
(30	2.5	1.75	0.5	1	0	7.5)

bowl pasta
(639.25	25.25	6.13	29.00	75.63	4.75	394.25)
1/4 cup mozzarella cheese
(80	5	3.5	6	1	0	190)
2 tbs parmesan cheese winco brand
(20	1.5	1	2	0	0	100)

cup of coffee cold brew with 1 tbs pure cane sugar
(45	0	0	0	12	0	0)
2 strawberry walmart brand poptarts of roommate's
(400	10	5	4	74	1	240)
=1152+120+200+106.7+330+45+30+639.3+80+20+45+400
=39.3+10+2+6.7+14+0+2.5+25.3+5+1.5+0+10
=10.0+7+0+4.7+11+0+1.8+6.1+3.5+1+0+5
=73.3+2+4+8+5+0+0.5+29+6+2+0+4
=137.44+4+42+1.3+44+12+1+75.6+1+0+12+74
=25.3+0+4+0+0+0+0+4.8+0+0+0+1
=1426.9+30+40+253.3+60+0+7.5+394.3+190+100+0+240



'@
$wsMeasures.Range("AA71").Value2 = $aaTextRow71

$wsMeasures.Range("AB71").Formula = "=1152+120+200+106.7+330+45+30+639.3+80+20+45+400"
$wsMeasures.Range("AC71").Formula = "=39.3+10+2+6.7+14+0+2.5+25.3+5+1.5+0+10"
$wsMeasures.Range("AD71").Formula = "=10+7+0+4.7+11+0+1.8+6.1+3.5+1+0+5"
$wsMeasures.Range("AE71").Formula = "=73.3+2+4+8+5+0+0.5+29+6+2+0+4"
$wsMeasures.Range("AF71").Formula = "=137.44+4+42+1.3+44+12+1+75.6+1+0+12+74"
$wsMeasures.Range("AG71").Formula = "=25.3+0+4+0+0+0+0+4.8+0+0+0+1"
$wsMeasures.Range("AH71").Formula = "=1426.9+30+40+253.3+60+0+7.5+394.3+190+100+0+240"

# --- researchMeasures: add new row 72 --------------------------------
# Copy formatting (styles + row height) from row 71 first
$wsMeasures.Rows.Item(71).Copy() | Out-Null
$wsMeasures.Rows.Item(72).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$zTextRow72 = @'
Woke up at 430 am to complete the worksheet on chemical reactions, stoichiometry, etc due today by 1159 pm. Got half way done. Mostly done in lab over last week but messy. Have to read through mess and calculate Molar mass and check items. Will finish the rest tonight. Had the rest of cold pressed about 1/2 cup with water, then 2 cups of my regular instant blend after making more of the cold pressed. Gave Mr. Growly his medicine when I woke up and the other when I get home around 5 pm and plan on working on homework till completed of this worksheet. I did my normal routine with feeding the babies after cleaning up some pet messes. Had a lg BM before finishing my 3rd cup of coffee for the day, then finished it. Checked the roommate's return he got and they took out an extra $100 approx bc they didn't charge his card even though the information was provided because the roommate wanted me to check it and asked why it was $100 less. Then took measurements around 705 am. then had a big bowl of spaghetti. 
'@
$aaTextRow72 = @'
bowl pasta 3-25-21 recipe
(639.25	25.25	6.13	29.00	75.63	4.75	394.25)
1/4 cup mozzarella cheese
(80	5	3.5	6	1	0	190)
85 degrees celsius sea salt coffee medium
(220    14       9      1	22	0	200)
chocolate hazelnut muffin tasted crappy, the oil was gross in the pan, tasted like mineral oil coated muffin, tried about 1/2 of it before throwing in the trash
(350	20.5	4	4	39	0	315)
hokkaido cheese tart
(230	18	10	5	20	0	140)
3/4 of 3 cheese quesadillas, bc tortillas broke off 1/4 of the shells in bag
 6*3/4=9/4 or 2 1/4 corn tortillas
(225	2.25	0	4.5	47.25	4.5	45)
1/2 cup mozzarella cheese
(160	10	7	12	2	0	380)
bowl pasta 3-25-21 recipe
(639.25	25.25	6.13	29.00	75.63	4.75	394.25)
=639.3+80+220+350+230+225+160+639.3
=25.3+5+14+20.5+18+2.25+10+25.3
=6.13+3.5+9+4+10+0+7+6.13
=29+6+1+4+5+4.5+12+29
=75.6+1+22+39+20+47.3+2+75.6
=4.75+0+0+0+0+4.5+0+4.75
=394.3+190+200+315+140+45+380+394.3

'@

$wsMeasures.Range("A72").Value2 = "Fri"
$wsMeasures.Range("B72").Value2 = 8
$wsMeasures.Range("C72").Value2 = 44281
$wsMeasures.Range("D72").Value2 = 0.2951388888888889
$wsMeasures.Range("E72").Value2 = 46
$wsMeasures.Range("F72").Value2 = 0
$wsMeasures.Range("G72").Value2 = 0
$wsMeasures.Range("H72").Value2 = 0
$wsMeasures.Range("I72").Value2 = 0
$wsMeasures.Range("J72").Value2 = 0.2951388888888889
$wsMeasures.Range("K72").Value2 = 137.19999999999999

$wsMeasures.Range("L72").Formula = "=K72-K71"
$wsMeasures.Range("M72").Formula = "=AB71"

$wsMeasures.Range("N72").Value2 = 30.75
$wsMeasures.Range("O72").Value2 = 32
$wsMeasures.Range("P72").Value2 = 10.5
$wsMeasures.Range("Q72").Value2 = 10.5
$wsMeasures.Range("R72").Value2 = 19
$wsMeasures.Range("S72").Value2 = 19.25
$wsMeasures.Range("T72").Value2 = 12
$wsMeasures.Range("U72").Value2 = 10
$wsMeasures.Range("V72").Value2 = 17
$wsMeasures.Range("W72").Value2 = 15
$wsMeasures.Range("X72").Value2 = 7
$wsMeasures.Range("Y72").Value2 = 7

$wsMeasures.Range("Z72").Value2 = $zTextRow72
$wsMeasures.Range("AA72").Value2 = $aaTextRow72

$wsMeasures.Range("AB72").Formula = "=639.3+80+220+350+230+225+160+639.3"
$wsMeasures.Range("AC72").Formula = "=25.3+5+14+20.5+18+2.25+10+25.3"
$wsMeasures.Range("AD72").Formula = "=6.13+3.5+9+4+10+0+7+6.13"
$wsMeasures.Range("AE72").Formula = "=29+6+1+4+5+4.5+12+29"
$wsMeasures.Range("AF72").Formula = "=75.6+1+22+39+20+47.3+2+75.6"
$wsMeasures.Range("AG72").Formula = "=4.75+0+0+0+0+4.5+0+4.75"
$wsMeasures.Range("AH72").Formula = "=394.3+190+200+315+140+45+380+394.3"

$wsMeasures.Range("AI72").Formula = "=`$AC72/`$AB72"
$wsMeasures.Range("AJ72").Formula = "=`$AD72/`$AB72"
$wsMeasures.Range("AK72").Formula = "=`$AE72/`$AB72"
$wsMeasures.Range("AL72").Formula = "=`$AF72/`$AB72"
$wsMeasures.Range("AM72").Formula = "=`$AG72/`$AB72"
$wsMeasures.Range("AN72").Formula = "=`$AH72/`$AB72"

$wsMeasures.Range("AO72").Value2 = 5
$wsMeasures.Range("AP72").Value2 = 1
$wsMeasures.Range("AQ72").Value2 = 1
$wsMeasures.Range("AR72").Value2 = 0
$wsMeasures.Range("AS72").Value2 = 0
$wsMeasures.Range("AT72").Value2 = 0
$wsMeasures.Range("AU72").Value2 = 0
$wsMeasures.Range("AV72").Value2 = 0
$wsMeasures.Range("AW72").Value2 = 31
$wsMeasures.Range("AX72").Value2 = 1
$wsMeasures.Range("AY72").Value2 = 7
$wsMeasures.Range("AZ72").Value2 = 0
$wsMeasures.Range("BA72").Value2 = 1
$wsMeasures.Range("BB72").Value2 = 0
$wsMeasures.Range("BC72").Value2 = 1
$wsMeasures.Range("BD72").Value2 = 1
$wsMeasures.Range("BE72").Value2 = 0
$wsMeasures.Range("BF72").Value2 = 0
$wsMeasures.Range("BG72").Value2 = 0
$wsMeasures.Range("BH72").Value2 = 0
$wsMeasures.Range("BI72").Value2 = 0

# --- NutritionalData: rewrite rows 187-189, add row 190 --------------
$wsNutrition.Range("A187").Value2 = '85 degrees celsius sea salt coffee medium'
$wsNutrition.Range("B187").Value2 = 220
$wsNutrition.Range("C187").Value2 = 14
$wsNutrition.Range("D187").Value2 = 9
$wsNutrition.Range("E187").Value2 = 1
$wsNutrition.Range("F187").Value2 = 22
$wsNutrition.Range("G187").Value2 = 0
$wsNutrition.Range("H187").Value2 = 200

$wsNutrition.Range("A188").Value2 = '85 degrees celsius hokkaido cheese tart'
$wsNutrition.Range("B188").Value2 = 230
$wsNutrition.Range("C188").Value2 = 18
$wsNutrition.Range("D188").Value2 = 10
$wsNutrition.Range("E188").Value2 = 5
$wsNutrition.Range("F188").Value2 = 20
$wsNutrition.Range("G188").Value2 = 0
$wsNutrition.Range("H188").Value2 = 140

$wsNutrition.Range("A189").Value2 = 'chocolate hazelnet muffin'
$wsNutrition.Range("B189").Value2 = 700
$wsNutrition.Range("C189").Value2 = 41
$wsNutrition.Range("D189").Value2 = 8
$wsNutrition.Range("E189").Value2 = 8
$wsNutrition.Range("F189").Value2 = 78
$wsNutrition.Range("G189").Formula = "=G187/2"
$wsNutrition.Range("H189").Value2 = 630

$wsNutrition.Range("B190").Formula = "=B57*2.25"
$wsNutrition.Range("C190").Formula = "=C57*2.25"
$wsNutrition.Range("D190").Formula = "=D57*2.25"
$wsNutrition.Range("E190").Formula = "=E57*2.25"
$wsNutrition.Range("F190").Formula = "=F57*2.25"
$wsNutrition.Range("G190").Formula = "=G57*2.25"
$wsNutrition.Range("H190").Formula = "=H57*2.25"

# --- View/selection adjustments ---------------------------------------
$wsNutrition.Activate()
$wsNutrition.Range("B190:H190").Select()

$wsMeasures.Activate()
$wsMeasures.Range("C72").Select()

